$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from column AL (38) to new column AM (39) for rows 1-27
$ws.Range("AL1:AL27").Copy()
$ws.Range("AM1:AM27").PasteSpecial(-4122)

# New date header for the added day (AM1)
$ws.Cells.Item(1, 39).Value = 45903

# Attendance values for the new date, one per player row
$ws.Cells.Item(2, 39).Value = "P"
$ws.Cells.Item(3, 39).Value = "P"
$ws.Cells.Item(4, 39).Value = "P"
$ws.Cells.Item(5, 39).Value = "P"
$ws.Cells.Item(6, 39).Value = "B"
$ws.Cells.Item(7, 39).Value = "P"
$ws.Cells.Item(8, 39).Value = "P"
$ws.Cells.Item(9, 39).Value = "P"
$ws.Cells.Item(10, 39).Value = "B"
$ws.Cells.Item(11, 39).Value = "P"
$ws.Cells.Item(12, 39).Value = "P"
$ws.Cells.Item(13, 39).Value = "P"
$ws.Cells.Item(14, 39).Value = "P"
$ws.Cells.Item(15, 39).Value = "P"
$ws.Cells.Item(16, 39).Value = "P"
$ws.Cells.Item(17, 39).Value = "RH"
$ws.Cells.Item(18, 39).Value = "P"
$ws.Cells.Item(19, 39).Value = "P"
$ws.Cells.Item(20, 39).Value = "P"
$ws.Cells.Item(21, 39).Value = "P"
$ws.Cells.Item(22, 39).Value = "P"
$ws.Cells.Item(23, 39).Value = "P"
$ws.Cells.Item(24, 39).Value = "P"
$ws.Cells.Item(25, 39).Value = "P"
$ws.Cells.Item(26, 39).Value = "P"
$ws.Cells.Item(27, 39).Value = "RH"

# Remove the old single-column total row (row 28) which is no longer needed
$ws.Rows.Item(28).Delete()

# Refresh the selection/active cell to match the latest edit location
$ws.Range("AP18").Select()
